# Remove duplicate Python events (rows 6-13 in the "events" sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:A13").EntireRow.Select()
$ws.Range("A6:A13").EntireRow.Delete()
